$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the ID_Employee value for the data row (A2), which was previously blank.
$ws.Range("A2").Value = "x007"
